$wb = $excel.ActiveWorkbook

$xlPasteValues = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# A reference sheet that already carries the "header row" / "index
# column" styling used throughout this workbook (bold, bordered,
# centered header; bordered index column) so the new sheet matches.
$refSheet = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------
# 2. Fill the new "2022-Q1" sheet with the same layout used by the other
#    quarterly fund-holding sheets (基金代码 / 基金名称 / 基金规模 / ...).
# ---------------------------------------------------------------------

# -- header row (B1:H1): copy the formatting used on the other sheets,
#    then fill in this sheet's own header text.
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# -- index cell A2: copy formatting from the reference sheet's A column.
$refSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial($xlPasteFormats)
$newSheet.Cells.Item(2, 1).Value = 0

# -- data row 2 (B2:G2): these columns are stored as plain text in this
#    workbook (even the numeric-looking ones), with no special cell
#    style. Stage the text in an out-of-the-way range formatted as
#    Text, then copy just the *values* over so the destination cells
#    pick up the text type without inheriting a new style.
$stage = $newSheet.Range("Z1:AE1")
$stage.NumberFormat = "@"
$stageValues = @("007497", "中庚价值灵动灵活配置混合", "24.35", "89.42", "2.21", "0.5381")
for ($i = 0; $i -lt $stageValues.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 26).Value = $stageValues[$i]
}
$stage.Copy()
$newSheet.Range("B2:G2").PasteSpecial($xlPasteValues)
$stage.Clear()

# -- H2 (仓位排名) is a plain number.
$newSheet.Cells.Item(2, 8).Value = 7

# ---------------------------------------------------------------------
# 3. Insert a new top data row in the "总计" sheet for the 2022-Q1 totals,
#    pushing the existing quarters down by one row.
#    (Re-fetch the sheet by name: after the insert above, the sheet's
#    tab position shifted, and stale worksheet handles track tab
#    position rather than identity.)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Excel's row-insert carries over neighbouring formatting onto the new
# row; this workbook's data rows (besides the index column) don't carry
# any explicit style, so strip whatever got inherited first.
$totalSheet.Range("A2:D2").ClearFormats()

# A2 needs the same bordered/centered style as the rest of the index
# column; copy it from the row that used to be A2 (now shifted to A3).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)
$totalSheet.Cells.Item(2, 1).Value = 0

# B2 (日期 column) is plain text, same treatment as above.
$stage2 = $totalSheet.Range("Z1")
$stage2.NumberFormat = "@"
$stage2.Value = "2022-Q1"
$stage2.Copy()
$totalSheet.Range("B2").PasteSpecial($xlPasteValues)
$stage2.Clear()

$totalSheet.Cells.Item(2, 3).Value = 1
$totalSheet.Cells.Item(2, 4).Value = 0.54

# Renumber the index column (A) for the rows that shifted down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
